# Update the cryptocurrency price/volume table (columns D "Price" and
# E "Volume(1h)") for rows 2-51 with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.199.17"
$ws.Range("E2").Value = "  -5.48%  "
$ws.Range("D3").Value = "3.343.41"
$ws.Range("E3").Value = "  -2.23%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'565.68"
$ws.Range("E5").Value = "  -2.34%  "
$ws.Range("D6").Value = "'130.86"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.344.66"
$ws.Range("E8").Value = "  -2.19%  "
$ws.Range("D9").Value = "'0.473"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").Value = "'7.44"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("E11").Value = "  -4.57%  "
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "3.911.88"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "3.342.13"
$ws.Range("E15").Value = "  -2.40%  "
$ws.Range("E16").Value = "  -3.58%  "
$ws.Range("D17").Value = "'24.71"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").Value = "60.259.63"
$ws.Range("E18").Value = "  -5.36%  "
$ws.Range("D19").Value = "'5.70"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").Value = "'13.48"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").Value = "'9.15"
$ws.Range("E21").Value = "  -6.68%  "
$ws.Range("D22").Value = "'355.42"
$ws.Range("E22").Value = "  -7.12%  "
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").Value = "3.476.45"
$ws.Range("E24").Value = "  -2.25%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "'69.42"
$ws.Range("E26").Value = "  -6.18%  "
$ws.Range("E27").Value = "  +2.66%  "
$ws.Range("D28").Value = "'1.69"
$ws.Range("E28").Value = "  +20.26%  "
$ws.Range("D29").Value = "'7.52"
$ws.Range("E29").Value = "  +7.48%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'7.96"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "3.374.66"
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("D37").Value = "'5.41"
$ws.Range("E37").Value = "  +5.71%  "
$ws.Range("D38").Value = "'6.91"
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("D40").Value = "'159.14"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").Value = "'0.0772"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "'4.39"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("E44").Value = "  +8.48%  "
$ws.Range("E45").Value = "  -4.18%  "
$ws.Range("D46").Value = "'40.83"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").Value = "'23.89"
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("D50").Value = "'22.43"
$ws.Range("E50").Value = "  +10.69%  "
$ws.Range("D51").Value = "'0.894"
$ws.Range("E51").Value = "  +0.86%  "
